$wb = $excel.ActiveWorkbook

# ALC row 4
$ws = $wb.Worksheets.Item(1)
$ws.Range("H4").Value = 1055.7142
$ws.Range("I4").Value = 680
$ws.Range("K4").Value = 680
$ws.Range("M4").Value = -566

# ALC row 18
$ws = $wb.Worksheets.Item(1)
$ws.Range("H18").Value = 1822
$ws.Range("I18").Value = 894.5
$ws.Range("K18").Value = 894.5
$ws.Range("M18").Value = -610.5

# ALC row 33
$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 146.07692
$ws.Range("I33").Value = 146.07692
$ws.Range("K33").Value = 146.07692
$ws.Range("M33").Value = 82.92308

# ALC row 86
$ws = $wb.Worksheets.Item(1)
$ws.Range("H86").Value = 252471.88
$ws.Range("I86").Value = 3217.75
$ws.Range("K86").Value = 3217.75
$ws.Range("M86").Value = -2094.75

# ALC row 89
$ws = $wb.Worksheets.Item(1)
$ws.Range("H89").Value = 252471.88
$ws.Range("I89").Value = 3217.75
$ws.Range("K89").Value = 16088.75
$ws.Range("M89").Value = -10472.75

# ALC row 111
$ws = $wb.Worksheets.Item(1)
$ws.Range("H111").Value = 49600
$ws.Range("I111").Value = 49500
$ws.Range("K111").Value = 148500
$ws.Range("M111").Value = -145433

# ALC row 113
$ws = $wb.Worksheets.Item(1)
$ws.Range("H113").Value = 4499
$ws.Range("I113").Value = 4499
$ws.Range("K113").Value = 4499
$ws.Range("M113").Value = -1245

# ALC row 135
$ws = $wb.Worksheets.Item(1)
$ws.Range("H135").Value = 1373.2667
$ws.Range("I135").Value = 468.7857
$ws.Range("K135").Value = 4219.071300000001
$ws.Range("M135").Value = -1684.071300000001

# ARM row 32
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 2335612.8
$ws.Range("I32").Value = 2189850.8
$ws.Range("K32").Value = 2189850.8
$ws.Range("M32").Value = -2189563.8

# ARM row 88
$ws = $wb.Worksheets.Item(2)
$ws.Range("H88").Value = 2423.4285
$ws.Range("J88").Value = 2243
$ws.Range("L88").Value = 2243
$ws.Range("N88").Value = -3055

# ARM row 91
$ws = $wb.Worksheets.Item(2)
$ws.Range("H91").Value = 2423.4285
$ws.Range("J91").Value = 2243
$ws.Range("L91").Value = 2243
$ws.Range("N91").Value = -5051

# ARM row 92
$ws = $wb.Worksheets.Item(2)
$ws.Range("H92").Value = 25000
$ws.Range("J92").Value = 25000
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992

# BSM row 99
$ws = $wb.Worksheets.Item(3)
$ws.Range("H99").Value = 1500
$ws.Range("I99").Value = 1500
$ws.Range("K99").Value = 1500
$ws.Range("M99").Value = -2

# BSM row 105
$ws = $wb.Worksheets.Item(3)
$ws.Range("H105").Value = 2875
$ws.Range("I105").Value = 2250
$ws.Range("K105").Value = 2250
$ws.Range("M105").Value = -503

# BSM row 107
$ws = $wb.Worksheets.Item(3)
$ws.Range("H107").Value = 299
$ws.Range("I107").Value = 299
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 299
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1621
$ws.Range("N107").ClearContents()

# CRP row 16
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 5541.3335
$ws.Range("I16").Value = 3312
$ws.Range("K16").Value = 3312
$ws.Range("M16").Value = -3025

# CRP row 41
$ws = $wb.Worksheets.Item(4)
$ws.Range("H41").Value = 1000
$ws.Range("I41").Value = 1000
$ws.Range("K41").Value = 1000
$ws.Range("M41").Value = -572

# CRP row 99
$ws = $wb.Worksheets.Item(4)
$ws.Range("H99").Value = 1347.8125
$ws.Range("I99").Value = 1251
$ws.Range("J99").Value = 2800
$ws.Range("K99").Value = 1251
$ws.Range("L99").Value = 2800
$ws.Range("M99").Value = 247
$ws.Range("N99").Value = -5796

# CRP row 105
$ws = $wb.Worksheets.Item(4)
$ws.Range("H105").Value = 2063.4546
$ws.Range("I105").Value = 1264
$ws.Range("J105").Value = 3218.2222
$ws.Range("K105").Value = 1264
$ws.Range("L105").Value = 3218.2222
$ws.Range("M105").Value = 483
$ws.Range("N105").Value = -6712.2222

# CRP row 108
$ws = $wb.Worksheets.Item(4)
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# CRP row 113
$ws = $wb.Worksheets.Item(4)
$ws.Range("H113").Value = 5541.3335
$ws.Range("I113").Value = 3312
$ws.Range("K113").Value = 3312
$ws.Range("M113").Value = -1142

# CRP row 115
$ws = $wb.Worksheets.Item(4)
$ws.Range("H115").Value = 44444
$ws.Range("J115").Value = 44444
$ws.Range("L115").Value = 44444
$ws.Range("N115").Value = -46794

# CRP row 126
$ws = $wb.Worksheets.Item(4)
$ws.Range("H126").Value = 1347.8125
$ws.Range("I126").Value = 1251
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 3753
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -1283
$ws.Range("N126").Value = -13340

# CUL row 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 572.5
$ws.Range("J5").Value = 572.5
$ws.Range("L5").Value = 1717.5
$ws.Range("N5").Value = -1941.5

# CUL row 52
$ws = $wb.Worksheets.Item(5)
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# CUL row 98
$ws = $wb.Worksheets.Item(5)
$ws.Range("H98").Value = 1559.8
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 1559.8
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 4679.4
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -7675.4

# CUL row 119
$ws = $wb.Worksheets.Item(5)
$ws.Range("H119").Value = 3332.3333
$ws.Range("I119").Value = 3332.3333
$ws.Range("K119").Value = 9996.999899999999
$ws.Range("M119").Value = -5158.999899999999

# CUL row 132
$ws = $wb.Worksheets.Item(5)
$ws.Range("H132").Value = 4499.6
$ws.Range("J132").Value = 6000
$ws.Range("L132").Value = 54000
$ws.Range("N132").Value = -59060

# CUL row 135
$ws = $wb.Worksheets.Item(5)
$ws.Range("H135").Value = 572.5
$ws.Range("J135").Value = 572.5
$ws.Range("L135").Value = 5152.5
$ws.Range("N135").Value = -10222.5

# GSM row 20
$ws = $wb.Worksheets.Item(6)
$ws.Range("H20").Value = 200000000
$ws.Range("I20").Value = 200000000
$ws.Range("K20").Value = 200000000
$ws.Range("M20").Value = -199999755

# GSM row 24
$ws = $wb.Worksheets.Item(6)
$ws.Range("H24").Value = 50011604
$ws.Range("I24").Value = 125007000
$ws.Range("J24").Value = 14673.667
$ws.Range("K24").Value = 125007000
$ws.Range("L24").Value = 14673.667
$ws.Range("M24").Value = -125006827
$ws.Range("N24").Value = -15019.667

# GSM row 44
$ws = $wb.Worksheets.Item(6)
$ws.Range("H44").Value = 26000
$ws.Range("I44").Value = 20000
$ws.Range("J44").Value = 50000
$ws.Range("K44").Value = 20000
$ws.Range("L44").Value = 50000
$ws.Range("M44").Value = -19404
$ws.Range("N44").Value = -51192

# GSM row 74
$ws = $wb.Worksheets.Item(6)
$ws.Range("H74").Value = 50000
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51872

# GSM row 77
$ws = $wb.Worksheets.Item(6)
$ws.Range("H77").Value = 50000
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -159360

# GSM row 97
$ws = $wb.Worksheets.Item(6)
$ws.Range("H97").Value = 280
$ws.Range("I97").Value = 280
$ws.Range("K97").Value = 280
$ws.Range("M97").Value = 216

# GSM row 102
$ws = $wb.Worksheets.Item(6)
$ws.Range("H102").Value = 2114.1538
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# GSM row 107
$ws = $wb.Worksheets.Item(6)
$ws.Range("H107").Value = 2315
$ws.Range("I107").Value = 763.6667
$ws.Range("J107").Value = 3866.3333
$ws.Range("K107").Value = 763.6667
$ws.Range("L107").Value = 3866.3333
$ws.Range("M107").Value = 1156.3333
$ws.Range("N107").Value = -7706.3333

# GSM row 113
$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 601.9091
$ws.Range("I113").Value = 490.625
$ws.Range("K113").Value = 490.625
$ws.Range("M113").Value = 1679.375

# LTW row 7
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 7622.2104
$ws.Range("J7").Value = 7603.5835
$ws.Range("L7").Value = 7603.5835
$ws.Range("N7").Value = -7827.5835

# LTW row 16
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 199.5
$ws.Range("I16").Value = 199.5
$ws.Range("K16").Value = 199.5
$ws.Range("M16").Value = -29.5

# LTW row 31
$ws = $wb.Worksheets.Item(7)
$ws.Range("H31").Value = 315
$ws.Range("I31").Value = 315
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 315
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -67
$ws.Range("N31").ClearContents()

# LTW row 41
$ws = $wb.Worksheets.Item(7)
$ws.Range("H41").Value = 300
$ws.Range("I41").Value = 300
$ws.Range("K41").Value = 300
$ws.Range("M41").Value = 138

# LTW row 46
$ws = $wb.Worksheets.Item(7)
$ws.Range("H46").Value = 1780.5
$ws.Range("I46").Value = 1174.5
$ws.Range("K46").Value = 1174.5
$ws.Range("M46").Value = -986.5

# LTW row 100
$ws = $wb.Worksheets.Item(7)
$ws.Range("H100").Value = 3249.5334
$ws.Range("I100").Value = 3191.25
$ws.Range("K100").Value = 3191.25
$ws.Range("M100").Value = -2650.25

# LTW row 122
$ws = $wb.Worksheets.Item(7)
$ws.Range("H122").Value = 5514.0625
$ws.Range("J122").Value = 7031.077
$ws.Range("L122").Value = 21093.231
$ws.Range("N122").Value = -25993.231

# LTW row 126
$ws = $wb.Worksheets.Item(7)
$ws.Range("H126").Value = 7622.2104
$ws.Range("J126").Value = 7603.5835
$ws.Range("L126").Value = 22810.7505
$ws.Range("N126").Value = -27750.7505

# WVR row 41
$ws = $wb.Worksheets.Item(8)
$ws.Range("H41").Value = 18517.5
$ws.Range("J41").Value = 18763
$ws.Range("L41").Value = 18763
$ws.Range("N41").Value = -19543

# WVR row 49
$ws = $wb.Worksheets.Item(8)
$ws.Range("H49").Value = 47499.5
$ws.Range("J49").Value = 47499.5
$ws.Range("L49").Value = 47499.5
$ws.Range("N49").Value = -47959.5

# WVR row 54
$ws = $wb.Worksheets.Item(8)
$ws.Range("H54").Value = 42999.8
$ws.Range("J54").Value = 42999.8
$ws.Range("L54").Value = 42999.8
$ws.Range("N54").Value = -44039.8

# WVR row 55
$ws = $wb.Worksheets.Item(8)
$ws.Range("H55").Value = 10761.75
$ws.Range("I55").Value = 4349
$ws.Range("J55").Value = 30000
$ws.Range("K55").Value = 4349
$ws.Range("L55").Value = 30000
$ws.Range("M55").Value = -4072
$ws.Range("N55").Value = -30554

# WVR row 62
$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 16628.715
$ws.Range("J62").Value = 14099.75
$ws.Range("L62").Value = 14099.75
$ws.Range("N62").Value = -15347.75

# WVR row 65
$ws = $wb.Worksheets.Item(8)
$ws.Range("H65").Value = 16628.715
$ws.Range("J65").Value = 14099.75
$ws.Range("L65").Value = 70498.75
$ws.Range("N65").Value = -76738.75
